$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in Start Date (E) / End Date (F) columns with dates pulled from the time clock ---
# Ordered row-by-row (row, StartDate, EndDate) so the write order always
# matches the task order top to bottom.
$dates = @(
  , @(3,  "Aug 10th 2020",       "October 15th 2020")
  , @(4,  "October 15th 2020",   "February 1st 2021")
  , @(5,  "November 15th 2020",  "December 1st 2020")
  , @(6,  "December 1st 2020",   "February 1st 2021")
  , @(7,  "January  15th 2021",  "February 1st 2021")
  , @(8,  "January  15th 2021",  "February 1st 2021")
  , @(9,  "Febrary 1st 2021",    "February 15th 2021")
  , @(10, "Febrary 1st 2021",    "February 15th 2021")
  , @(11, "Febrary 1st 2021",    "March 1st 2021")
  , @(12, "Febrary 1st 2021",    "March 1st 2021")
  , @(13, "March 1st 2021",      "March 15th 2021")
  , @(14, "March 15th 2021",     "March 30th 2021")
  , @(15, "April 1st 2021",      "April 11th 2021")
)

foreach ($entry in $dates) {
    $row = $entry[0]
    $ws.Cells.Item($row, 5).Value = $entry[1]
    $ws.Cells.Item($row, 6).Value = $entry[2]
}

# Apply the custom date-style number format (text values display as-is, but
# the cells now carry a date format like the rest of the workbook's refresh).
# (covers the header row too, matching the header cells picking up the format)
$ws.Range("E2:F15").NumberFormat = "[$-409]mmmm\ d\,\ yyyy;@"

# Widen the Start/End Date columns now that they hold real content.
$ws.Columns("E:F").ColumnWidth = 18.6

# Row 9 loses its manual 30pt height now that the row has been touched again.
$ws.Rows(9).AutoFit() | Out-Null

# --- Window/view refresh recorded by Excel on save ---
$excel.ActiveWindow.Zoom = 70
$ws.Range("J7").Select() | Out-Null
